$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.131.18"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.478.40"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "584.33"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "172.87"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "2.478.22"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "2.930.78"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "25.49"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "67.092.19"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "2.467.26"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "10.96"
$ws.Range("E20").Value = "  -2.87%  "
$ws.Range("D21").Value = "350.04"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "68.92"
$ws.Range("D25").Value = "4.22"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "2.604.87"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "504.20"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "162.34"
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "18.15"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.69"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "0.328"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").Value = "2.37"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").Value = "143.09"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").Value = "0.0₆0264"
$ws.Range("E47").Value = "  +5.22%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "0.514"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("E51").Value = "  -1.10%  "
